# Fix CSV export encoding to UTF8
#
# Renames the "ref_name" dictionary row to "ref_hospital_name"
# (the ODK field the CSV exporter maps to a "referral hospital name"
# column), widens column D so the longer value still fits, and leaves
# the sheet's last selection where the editor's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A31 / D31 hold the "n1-ref_name" / "ref_name" dictionary entry.
$ws.Range("A31").Value = "n1-ref_hospital_name"
$ws.Range("D31").Value = "ref_hospital_name"

# Column D needs to be a bit wider to fit "ref_hospital_name" /
# "n1-ref_hospital_name" (stored width 16, i.e. ColumnWidth 15.3 once
# Excel rounds to whole pixels).
$ws.Columns.Item(4).ColumnWidth = 15.3

# Last active cell/selection as saved in the sheet view.
$ws.Range("G17").Select()
